$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-7 from 45233 to 45243
$ws.Range("C2:C7").Value = 45243
